# Replicates the "Added code from slides for initial putexcel example as a
# do-file" commit: the regression output was re-run (Stata putexcel), which
# overwrote the row-name labels on "Regression Table 1" and the generated
# timestamp footers on the two "Regression Shell" tables, and left the
# first table's window/tab active (with cell P7 selected) instead of the
# third table.

$wb = $excel.ActiveWorkbook

# --- "Regression Shell Table 1" / "Regression Shell Table 2": refresh the
# putexcel-generated timestamp footer in A30. Do this before touching the
# first sheet so the shared-string table lands in the same append order as
# the source commit (dates first, then the new row-name labels).
$wsShell1 = $wb.Worksheets.Item(2)
$wsShell1.Range("A30").Value2 = "08:18:46 15 Nov 2018"

$wsShell2 = $wb.Worksheets.Item(3)
$wsShell2.Range("A30").Value2 = "08:19:59 15 Nov 2018"

# --- "Regression Table 1": the raw regression dump now carries Stata's own
# coefficient names ("smoke", "_cons") instead of the generic "rownames(temp)"
# placeholder that was in every row before.
$wsReg = $wb.Worksheets.Item(1)
$wsReg.Range("A5").Value2 = "smoke"
$wsReg.Range("A6").Value2 = "_cons"

# --- Active window/tab moved to "Regression Table 1", with the selection
# left at P7.
$wsReg.Activate()
$wsReg.Range("P7").Select()
